$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '26.888.81'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  +0.12%  '
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '1.547.41'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  -1.05%  '
$ws.Range("E4").Value = '  +0.07%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '206.36'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +0.08%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '0.487'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -0.35%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("E8").Value = '  -0.25%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '21.50'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  -1.83%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.0583'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  -0.38%  '
$ws.Range("E11").Value = '  -1.50%  '
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '1.767.70'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  -1.08%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '1.545.27'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  -1.23%  '
$ws.Range("E14").Value = '  -1.33%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '0.512'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -0.25%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '26.902.24'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  +0.12%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '61.45'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  +0.25%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '215.06'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +0.21%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '7.25'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  -1.73%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '0.0{0}0683' -f [char]0x2083
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +0.83%  '
$ws.Range("E21").Value = '  +0.13%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '4.02'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -2.81%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '9.18'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -1.03%  '
$ws.Range("E24").Value = '  -2.96%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '152.68'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  -0.58%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '6.65'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  -0.93%  '
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '14.87'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -0.59%  '
$ws.Range("E28").Value = '  +0.06%  '
$ws.Range("E29").Value = '  -0.38%  '
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '0.0460'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  -0.97%  '
$ws.Range("E31").Value = '  -1.16%  '
$ws.Range("E32").Value = '  +1.68%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '1.367.32'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -2.56%  '
$ws.Range("E34").Value = '  +0.57%  '
$ws.Range("E35").Value = '  -0.25%  '
$ws.Range("E36").Value = '  +2.83%  '
$ws.Range("E37").Value = '  -0.02%  '
$ws.Range("E38").Value = '  +0.87%  '
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '0.520'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  -1.82%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.809'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  -0.54%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '5.62'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +5.90%  '
$ws.Range("B42").Value = 'WEMIXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.991'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +0.32%  '
$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '2.22'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  +1.96%  '
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '63.53'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +0.75%  '
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '1.72'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  -4.23%  '
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '1.682.42'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -1.16%  '
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '84.74'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -1.98%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '0.0505'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  +2.72%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '0.0{0}0965' -f [char]0x2087
$cell.Style = "Normal"
$ws.Range("E49").Value = '  -1.78%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '0.0947'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -0.14%  '
$ws.Range("E51").Value = '  -0.04%  '
